# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-02 (serial 45171) to 2023-09-03 (serial 45172).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45171) {
        $cell.Value = 45172
    }
}
